# Apply data updates to the "Inscricoes" sheet as described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 9
$ws.Range("E9").Value = 28
$ws.Range("F9").Value = 14
$ws.Range("H9").Value = 23

# Row 15
$ws.Range("E15").Value = 166

# Row 18
$ws.Range("E18").Value = 117

# Row 19
$ws.Range("E19").Value = 62

# Row 21
$ws.Range("E21").Value = 2

# Row 27
$ws.Range("E27").Value = 10

# Row 32
$ws.Range("E32").Value = 21

# Row 33
$ws.Range("E33").Value = 46

# Row 36
$ws.Range("E36").Value = 110
$ws.Range("F36").Value = 50
$ws.Range("H36").Value = 82

# Row 37
$ws.Range("E37").Value = 58

# Row 52
$ws.Range("E52").Value = 8
$ws.Range("F52").Value = 4
$ws.Range("H52").Value = 4

# Row 64
$ws.Range("E64").Value = 35

# Row 69
$ws.Range("E69").Value = 16

# Row 77
$ws.Range("E77").Value = 56

# Row 78
$ws.Range("E78").Value = 47
$ws.Range("F78").Value = 22
$ws.Range("H78").Value = 43

# Row 84
$ws.Range("E84").Value = 5

# Row 88
$ws.Range("E88").Value = 27
$ws.Range("F88").Value = 15
$ws.Range("H88").Value = 23
